# Update "Algorithms tests comparison table" worksheet:
#  - Drop the now-unused "Right hand punch"/"Left hand punch"/"Jump" header
#    columns and the "Minimal-Variance matching" row; add the
#    "Taekwondo position" column with measured results for the two
#    remaining algorithms (Dynamic time warping / Elastic action
#    comparison with freedom degree), and bold the header/label cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 ("Minimal-Variance matching") is removed entirely; this shifts
# row 4 ("Elastic action comparison with freedom degree") up to row 3.
$ws.Rows(3).Delete()

# --- Results for "Dynamic time warping" (row 2) ---
$ws.Range("B2").Value = 0.866
$ws.Range("B2").NumberFormat = "0.00%"
$ws.Range("C2").Value = 1
$ws.Range("C2").NumberFormat = "0%"
$ws.Range("D2").Value = 1
$ws.Range("D2").NumberFormat = "0%"
$ws.Range("E2").Value = 0.9166
$ws.Range("E2").NumberFormat = "0.00%"
$ws.Range("F2").Value = 0.7143
$ws.Range("F2").NumberFormat = "0.00%"

# --- Results for "Elastic action comparison with freedom degree" (row 3) ---
$ws.Range("B3").Value = 0.8
$ws.Range("B3").NumberFormat = "0.00%"
$ws.Range("C3").Value = 0.75
$ws.Range("C3").NumberFormat = "0.00%"
$ws.Range("D3").Value = 0.8182
$ws.Range("D3").NumberFormat = "0.00%"
$ws.Range("E3").Value = 1
$ws.Range("E3").NumberFormat = "0.00%"
$ws.Range("F3").Value = 0.7557
$ws.Range("F3").NumberFormat = "0.00%"

# --- Header row: bold it, swap the trailing three test columns for a
# single "Taekwondo position" column, and clear the now-empty H1/I1 ---
$ws.Range("A1").Font.Bold = $true
$ws.Range("B1:G1").Font.Bold = $true
$ws.Range("G1").Value = "Taekwondo position"
$ws.Range("H1").ClearContents()
$ws.Range("I1").ClearContents()

# --- Bold the algorithm-name labels in column A ---
$ws.Range("A2").Font.Bold = $true
$ws.Range("A3").Font.Bold = $true

# Slightly widen column B and leave the selection on H1, matching the
# saved view state.
$ws.Columns("B").ColumnWidth = 6.666666666666667
$ws.Range("H1").Select() | Out-Null
